$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# --- Copy formatting (incl. the thin-border style used on column D) from the
# last existing pair of rows (198:199) down across the 12 new rows (200:211)
# so the new rows pick up the already-existing style index instead of Excel
# minting a brand new one. ---
$ws.Range("A198:E199").Copy()
$ws.Range("A200:E211").PasteSpecial(-4122)

# --- New product block: "Previdência" under "Seguros" ---

# Icatu (até R$299,99)
$ws.Range("A200").Value = "Seguros"
$ws.Range("B200").Value = "Previdência"
$ws.Range("C200").Value = "Icatu (até R$299,99)"
$ws.Range("D200").Value = 1
$ws.Range("E200").Value = 0.1

$ws.Range("A201").Value = "Seguros"
$ws.Range("B201").Value = "Previdência"
$ws.Range("C201").Value = "Icatu (até R$299,99)"
$ws.Range("D201").Value = 2
$ws.Range("E201").Value = 0

# Icatu (R$300,00 - R$599,99)
$ws.Range("A202").Value = "Seguros"
$ws.Range("B202").Value = "Previdência"
$ws.Range("C202").Value = "Icatu (R$300,00 - R$599,99)"
$ws.Range("D202").Value = 1
$ws.Range("E202").Value = 0.25

$ws.Range("A203").Value = "Seguros"
$ws.Range("B203").Value = "Previdência"
$ws.Range("C203").Value = "Icatu (R$300,00 - R$599,99)"
$ws.Range("D203").Value = 2
$ws.Range("E203").Value = 0

# Icatu (apartir de R$600,00)
$ws.Range("A204").Value = "Seguros"
$ws.Range("B204").Value = "Previdência"
$ws.Range("C204").Value = "Icatu (apartir de R$600,00)"
$ws.Range("D204").Value = 1
$ws.Range("E204").Value = 0.6

$ws.Range("A205").Value = "Seguros"
$ws.Range("B205").Value = "Previdência"
$ws.Range("C205").Value = "Icatu (apartir de R$600,00)"
$ws.Range("D205").Value = 2
$ws.Range("E205").Value = 0

# Icatu Esporádico
$ws.Range("A206").Value = "Seguros"
$ws.Range("B206").Value = "Previdência"
$ws.Range("C206").Value = "Icatu Esporádico"
$ws.Range("D206").Value = 1
$ws.Range("E206").Value = 0.007

$ws.Range("A207").Value = "Seguros"
$ws.Range("B207").Value = "Previdência"
$ws.Range("C207").Value = "Icatu Esporádico"
$ws.Range("D207").Value = 2
$ws.Range("E207").Value = 0

# Sulamérica Prestige (até R$5000,00)
$ws.Range("A208").Value = "Seguros"
$ws.Range("B208").Value = "Previdência"
$ws.Range("C208").Value = "Sulamérica Prestige (até R$5000,00)"
$ws.Range("D208").Value = 1
$ws.Range("E208").Value = 1

$ws.Range("A209").Value = "Seguros"
$ws.Range("B209").Value = "Previdência"
$ws.Range("C209").Value = "Sulamérica Prestige (até R$5000,00)"
$ws.Range("D209").Value = 2
$ws.Range("E209").Value = 0

# Sulamérica Prestige Esporádico
$ws.Range("A210").Value = "Seguros"
$ws.Range("B210").Value = "Previdência"
$ws.Range("C210").Value = "Sulamérica Prestige Esporádico"
$ws.Range("D210").Value = 1
$ws.Range("E210").Value = 0.015

$ws.Range("A211").Value = "Seguros"
$ws.Range("B211").Value = "Previdência"
$ws.Range("C211").Value = "Sulamérica Prestige Esporádico"
$ws.Range("D211").Value = 2
$ws.Range("E211").Value = 0

# --- Column C widened to fit the longest new label ("Sulamérica Prestige
# (até R$5000,00)") ---
$ws.Columns.Item(3).ColumnWidth = 32.166666666666664

# --- Update the view: scroll target + active selection mirror what the
# author had on screen after appending the rows ---
$ws.Range("E211").Select()

Write-Host "done"
